$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (before old row 2)
$ws.Rows("2:4").Insert()
# The insert copies formatting from the row above (the bold header row);
# clear it so the new data rows are unstyled, matching the rest of the data.
$ws.Range("A2:C4").ClearFormats()

# Populate the newly inserted rows with the new data points
$ws.Cells.Item(2, 1).Value2 = -0.01334476470947288
$ws.Cells.Item(2, 2).Value2 = 0.1471533775329589
$ws.Cells.Item(2, 3).Value2 = -0.2691573500633239
$ws.Cells.Item(3, 1).Value2 = -0.2625431418418885
$ws.Cells.Item(3, 2).Value2 = 0.13736093044281
$ws.Cells.Item(3, 3).Value2 = -0.1842701695859431
$ws.Cells.Item(4, 1).Value2 = -0.3044750690460205
$ws.Cells.Item(4, 2).Value2 = 0.2236802577972413
$ws.Cells.Item(4, 3).Value2 = -0.08151795715093602

# Append 7 new rows of data after the existing data (old row 21, now row 24)
$ws.Cells.Item(25, 1).Value2 = 0.09789943695068341
$ws.Cells.Item(25, 2).Value2 = 0.3452561050653447
$ws.Cells.Item(25, 3).Value2 = -0.131537172943354
$ws.Cells.Item(26, 1).Value2 = 0.05506801605224545
$ws.Cells.Item(26, 2).Value2 = 0.07766664028167738
$ws.Cells.Item(26, 3).Value2 = -0.2464380264282224
$ws.Cells.Item(27, 1).Value2 = -0.06795549392700258
$ws.Cells.Item(27, 2).Value2 = 0.1180151626467706
$ws.Cells.Item(27, 3).Value2 = -0.2082828953862188
$ws.Cells.Item(28, 1).Value2 = -0.01404476165771439
$ws.Cells.Item(28, 2).Value2 = 0.2834141030907641
$ws.Cells.Item(28, 3).Value2 = -0.1360972765833135
$ws.Cells.Item(29, 1).Value2 = -0.01649236679077155
$ws.Cells.Item(29, 2).Value2 = 0.2205449156463143
$ws.Cells.Item(29, 3).Value2 = -0.119759158231318
$ws.Cells.Item(30, 1).Value2 = -0.003359794616699139
$ws.Cells.Item(30, 2).Value2 = 0.1425043791532516
$ws.Cells.Item(30, 3).Value2 = -0.1422623544931412
$ws.Cells.Item(31, 1).Value2 = 0.02320241928100578
$ws.Cells.Item(31, 2).Value2 = 0.1721755955368283
$ws.Cells.Item(31, 3).Value2 = -0.1290906090289352
